$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Shreyas Gopal"

# Insert a new column before column A for matchNo, shifting existing
# columns (teamName..result) one to the right (A:L -> B:M)
$ws.Columns.Item(1).Insert()

# Insert a new row before row 2 so the existing data row (now B2:M2)
# shifts down to B3:M3, leaving row 2 free for the new match record
$ws.Rows.Item(2).Insert()

# Header row
$ws.Range("A1").Value = "matchNo"

# Row 2 - new match record ("51st" innings vs Mumbai Indians)
$ws.Range("A2").Value = "51st"
$ws.Range("B2").Value = "Rajasthan Royals"
$ws.Range("C2").Value = "Shreyas Gopal"
$ws.Range("D2").Value = "c †Ishan Kishan b Bumrah"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "0"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "0.00"
$ws.Range("J2").Value = "Mumbai Indians"
$ws.Range("K2").Value = "Sharjah"
$ws.Range("L2").Value = "October 05"
$ws.Range("M2").Value = "Mumbai won by 8 wickets (with 70 balls remaining)"

# Row 3 - pre-existing match record ("16th" innings vs Royal Challengers
# Bangalore), shifted down from row 2; only the new matchNo cell is added
$ws.Range("A3").Value = "16th"
